# Refresh cached Market Board pricing + recomputed Leve profit columns (H:N)
# across the per-class Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Mirrors the scheduled market-data-refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising / Quicksilver
$ws.Range("H2").Value = 533.3333
$ws.Range("J2").Value = 1325
$ws.Range("L2").Value = 1325
$ws.Range("N2").Value = -1551
# Row 4: Root Rush / Growth Formula Alpha
$ws.Range("H4").Value = 57.75
$ws.Range("I4").Value = 57.75
$ws.Range("K4").Value = 57.75
$ws.Range("M4").Value = 56.25
# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 2333.3333
$ws.Range("I18").Value = 2000
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -1716
$ws.Range("N18").Value = -3568
# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 4725
$ws.Range("I43").Value = 5850
$ws.Range("J43").Value = 2475
$ws.Range("K43").Value = 5850
$ws.Range("L43").Value = 2475
$ws.Range("M43").Value = -5781
$ws.Range("N43").Value = -2613
# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 4624.5835
$ws.Range("I64").Value = 4499.4443
$ws.Range("K64").Value = 4499.4443
$ws.Range("M64").Value = -4251.4443
# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 4624.5835
$ws.Range("I67").Value = 4499.4443
$ws.Range("K67").Value = 4499.4443
$ws.Range("M67").Value = -3641.4443
# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 18459954
$ws.Range("J88").Value = 3257980
$ws.Range("L88").Value = 3257980
$ws.Range("N88").Value = -3258792
# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 18459954
$ws.Range("J91").Value = 3257980
$ws.Range("L91").Value = 3257980
$ws.Range("N91").Value = -3260788
# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 2387.5
$ws.Range("I106").Value = 2509.182
$ws.Range("J106").Value = 2119.8
$ws.Range("K106").Value = 2509.182
$ws.Range("L106").Value = 2119.8
$ws.Range("M106").Value = -1878.182
$ws.Range("N106").Value = -3381.8
# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 25000390
$ws.Range("I135").Value = 26316182
$ws.Range("K135").Value = 236845638
$ws.Range("M135").Value = -236843103
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2045.5
$ws.Range("I138").Value = 1718.579
$ws.Range("K138").Value = 5155.737
$ws.Range("M138").Value = -15.73700000000008
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 3724.2
$ws.Range("I141").Value = 3724.2
$ws.Range("K141").Value = 11172.6
$ws.Range("M141").Value = -5992.599999999999
$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 2150.5
$ws.Range("I5").Value = 1118.1666
$ws.Range("K5").Value = 1118.1666
$ws.Range("M5").Value = -1006.1666
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 38463940
$ws.Range("I61").Value = 41668268
$ws.Range("K61").Value = 41668268
$ws.Range("M61").Value = -41668056
# Row 80: A Squire to Inspire / Titanium Hoplon
$ws.Range("H80").Value = 49990
$ws.Range("J80").Value = 49990
$ws.Range("L80").Value = 49990
$ws.Range("N80").Value = -51986
# Row 83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws.Range("H83").Value = 49990
$ws.Range("J83").Value = 49990
$ws.Range("L83").Value = 149970
$ws.Range("N83").Value = -159954
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 38463940
$ws.Range("I136").Value = 41668268
$ws.Range("K136").Value = 125004804
$ws.Range("M136").Value = -125002254
$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 2150.5
$ws.Range("I4").Value = 1118.1666
$ws.Range("K4").Value = 1118.1666
$ws.Range("M4").Value = -1003.1666
# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Range("H82").Value = 11551.8
$ws.Range("I82").Value = 11551.8
$ws.Range("K82").Value = 11551.8
$ws.Range("M82").Value = -11168.8
# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Range("H85").Value = 11551.8
$ws.Range("I85").Value = 11551.8
$ws.Range("K85").Value = 11551.8
$ws.Range("M85").Value = -10225.8
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2450.8333
$ws.Range("I86").Value = 2219.182
$ws.Range("J86").Value = 4999
$ws.Range("K86").Value = 2219.182
$ws.Range("L86").Value = 4999
$ws.Range("M86").Value = -1096.182
$ws.Range("N86").Value = -7245
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2450.8333
$ws.Range("I89").Value = 2219.182
$ws.Range("J89").Value = 4999
$ws.Range("K89").Value = 11095.91
$ws.Range("L89").Value = 24995
$ws.Range("M89").Value = -5479.91
$ws.Range("N89").Value = -36227
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 14709450
$ws.Range("I134").Value = 16132600
$ws.Range("J134").Value = 3566.3333
$ws.Range("K134").Value = 48397800
$ws.Range("L134").Value = 10698.9999
$ws.Range("M134").Value = -48395265
$ws.Range("N134").Value = -15768.9999
$ws = $wb.Worksheets.Item("CRP")
# Row 5: Bowing Out / Maple Shortbow
$ws.Range("H5").Value = 5650.125
$ws.Range("I5").Value = 2249.1667
$ws.Range("J5").Value = 15853
$ws.Range("K5").Value = 2249.1667
$ws.Range("L5").Value = 15853
$ws.Range("M5").Value = -2137.1667
$ws.Range("N5").Value = -16077
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 1275.7693
$ws.Range("I7").Value = 1809.5555
$ws.Range("K7").Value = 1809.5555
$ws.Range("M7").Value = -1696.5555
# Row 19: Shielding Sales / Square Ash Shield
$ws.Range("H19").Value = 1580.25
$ws.Range("I19").Value = 1489.4
$ws.Range("J19").Value = 1731.6666
$ws.Range("K19").Value = 1489.4
$ws.Range("L19").Value = 1731.6666
$ws.Range("M19").Value = -1319.4
$ws.Range("N19").Value = -2071.6666
# Row 24: What You Need / Square Ash Shield
$ws.Range("H24").Value = 1580.25
$ws.Range("I24").Value = 1489.4
$ws.Range("J24").Value = 1731.6666
$ws.Range("K24").Value = 1489.4
$ws.Range("L24").Value = 1731.6666
$ws.Range("M24").Value = -1319.4
$ws.Range("N24").Value = -2071.6666
# Row 29: Grinding It Out / Mudstone Grinding Wheel
$ws.Range("H29").Value = 3757.5
$ws.Range("I29").Value = 15
$ws.Range("K29").Value = 15
$ws.Range("M29").Value = 278
# Row 59: Bow Down to Magic / Crab Bow
$ws.Range("H59").Value = 81333
$ws.Range("I59").Value = 4000
$ws.Range("J59").Value = 119999.5
$ws.Range("K59").Value = 4000
$ws.Range("L59").Value = 119999.5
$ws.Range("M59").Value = -2855
$ws.Range("N59").Value = -122289.5
# Row 68: Do You Even String Bow / Holy Cedar Composite Bow
$ws.Range("H68").Value = 98245.836
$ws.Range("J68").Value = 98245.836
$ws.Range("L68").Value = 98245.836
$ws.Range("N68").Value = -99743.836
# Row 71: Win One Bow, Get Three Free (L) / Holy Cedar Composite Bow
$ws.Range("H71").Value = 98245.836
$ws.Range("J71").Value = 98245.836
$ws.Range("L71").Value = 294737.508
$ws.Range("N71").Value = -302225.508
# Row 74: License to Heal / Dark Chestnut Rod
$ws.Range("H74").Value = 43915
$ws.Range("J74").Value = 55675
$ws.Range("L74").Value = 55675
$ws.Range("N74").Value = -57423
# Row 77: Purified Polyrhythm (L) / Dark Chestnut Rod
$ws.Range("H77").Value = 43915
$ws.Range("J77").Value = 55675
$ws.Range("L77").Value = 167025
$ws.Range("N77").Value = -175761
$ws = $wb.Worksheets.Item("CUL")
# Row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws.Range("H138").Value = 1598.75
$ws.Range("I138").Value = 1598.75
$ws.Range("K138").Value = 4796.25
$ws.Range("M138").Value = 343.75
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 2800
$ws.Range("J80").Value = 2800
$ws.Range("L80").Value = 2800
$ws.Range("N80").Value = -4796
# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 2800
$ws.Range("J83").Value = 2800
$ws.Range("L83").Value = 14000
$ws.Range("N83").Value = -23984
# Row 131: Star Athletes / Star Quartz Wristband of Aiming
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 701
$ws.Range("K55").Value = 701
$ws.Range("M55").Value = -528
$ws = $wb.Worksheets.Item("WVR")
# Row 69: Fashion Patrol / Holy Rainbow Sarouel of Casting
$ws.Range("H69").Value = 26666
$ws.Range("J69").Value = 26666
$ws.Range("L69").Value = 26666
$ws.Range("N69").Value = -28164
# Row 72: Dress Code Violation (L) / Holy Rainbow Sarouel of Casting
$ws.Range("H72").Value = 26666
$ws.Range("J72").Value = 26666
$ws.Range("L72").Value = 79998
$ws.Range("N72").Value = -87486
